$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: E4 date (was shared-string "14/10/2025") -> numeric date serial for 15/10/2025 ---
$ws.Range("E4").Value = 45945

# --- Row 13: F13 start time 08:00 -> 12:00 ---
$ws.Range("F13").Value = 0.5

# --- Row 28: F28 start time 02:00 -> 14:00 ---
$ws.Range("F28").Value = 0.58333333333333337

# --- Row 89: E89 date (was shared-string "05/10/2025") -> numeric date serial for 08/10/2025 ---
$ws.Range("E89").Value = 45938

# --- Row 90: E90 date (was shared-string "30/10/2025") -> numeric date serial for 28/10/2025 ---
$ws.Range("E90").Value = 45958
# F90 also needs to pick up F89's (white/no-fill) number style before its value changes
$ws.Range("F89").Copy()
$ws.Range("F90").PasteSpecial(-4122)
$ws.Range("F90").Value = 0.41666666666666669

# --- Row 91: E91 date (was shared-string "08/10/2025") -> numeric date serial for 29/09/2025 ---
$ws.Range("E91").Value = 45929
$ws.Range("F91").Value = 0.58333333333333337

# --- Row 121: F121 start time 02:00 -> 14:00 ---
$ws.Range("F121").Value = 0.58333333333333337

# --- Row 137: E137 date text changes from "03/11/2025" to "05/11/2025" (stays a shared string) ---
$ws.Range("E11").Copy()
$ws.Range("E137").PasteSpecial(-4163)

# --- Row 176: E176 date text changes from "29/10/2025" to "28/10/2025" and picks up the
#     no-fill date style (style of E89); F176 start time 14:00 -> 08:00 ---
$ws.Range("E89").Copy()
$ws.Range("E176").PasteSpecial(-4122)
$ws.Range("E176").Value = "28/10/2025"
$ws.Range("F176").Value = 0.33333333333333331

# --- Row 185: E185 date (was shared-string "16/10/2025") -> numeric date serial for 13/10/2025 ---
$ws.Range("E185").Value = 45943
$ws.Range("F185").Value = 0.375

$excel.CutCopyMode = $false

# --- Apply the Sort that was recorded for the table (data stays in the same order) ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B187"))
$ws.Sort.SortFields.Add($ws.Range("C2:C187"))
$ws.Sort.SortFields.Add($ws.Range("D2:D187"))
$ws.Sort.SortFields.Add($ws.Range("E2:E187"))
$ws.Sort.SetRange($ws.Range("A1:G187"))
$ws.Sort.Header = 1
$ws.Sort.Apply()
